$wb = $excel.ActiveWorkbook

# --- Add the new "table_mitt" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "table_mitt"

# --- Header row ---
$newSheet.Range("A1").Value = "LUC_Ag"
$newSheet.Range("B1").Value = "LUC_Crp"

# --- Descriptive text row ---
$newSheet.Range("A2").Value = "Carbon taxes are applied on emissions from deforestation. Carbon taxes are also applied to the emissions from conversion of natural land to cropland or grassland in Zambia. Emissions from converting grassland to cropland are not taxed. This scenario reflects a holistic policy on LUC targeting all agricultural area. "
$newSheet.Range("B2").Value = "Carbon taxes are applied on emissions from deforestation. Carbon taxes are also applied to emissions from the conversion of grassland to cropland and from natural land to cropland in Zambia. Emissions from the conversion of natural land to grassland are not taxed. This scenario reflects a targeted policy on LUC from only cropland expansion."

# --- Formatting: wrap text on the description row, taller row to fit the text ---
$newSheet.Range("A2:B2").WrapText = $true
$newSheet.Rows.Item(2).RowHeight = 105

# --- Column widths ---
$newSheet.Columns.Item(1).ColumnWidth = 44.5
$newSheet.Columns.Item(2).ColumnWidth = 62.83333333333333

# --- Selection / view state matches the source sheet ---
$newSheet.Range("A2").Select() | Out-Null

# Make the new sheet the active tab (moves tabSelected from table_3 to table_mitt)
$newSheet.Activate()

# --- table_3 gained an explicit width for column A ---
$tbl3 = $wb.Worksheets.Item("table_3")
$tbl3.Columns.Item(1).ColumnWidth = 9.7109375
